# Add new columns I (I0) and J (IF) to Sheet1, mirroring the existing
# header style used by columns B..H, and fill in the per-row data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
$ws.Cells.Item(1, 9).Value  = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Copy the formatting (bold, centered, bordered) from the existing H1
# header cell onto the two new header cells so they match the rest of
# the header row.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data ---------------------------------------------------------------
# Each entry is (row, I-value, J-value)
$data = @(
    @(2, 7, 7),
    @(3, 6, 6),
    @(4, 7, 7),
    @(5, 8, 8),
    @(6, 8, 8),
    @(7, 7, 8),
    @(8, 6, 7),
    @(9, 7, 7),
    @(10, 7, 7),
    @(11, 7, 7),
    @(12, 8, 8),
    @(13, 7, 7),
    @(14, 1, 1),
    @(15, 7, 7),
    @(16, 3, 5),
    @(17, 4, 5),
    @(18, 5, 5),
    @(19, 8, 8),
    @(20, 6, 6),
    @(21, 5, 6),
    @(22, 7, 7),
    @(23, 7, 7),
    @(24, 8, 8),
    @(25, 8, 8),
    @(26, 7, 7),
    @(27, 4, 5),
    @(28, 8, 8),
    @(29, 7, 7),
    @(30, 6, 6),
    @(31, 3, 4),
    @(32, 5, 5),
    @(33, 8, 8),
    @(34, 7, 7),
    @(35, 7, 7),
    @(36, 7, 7),
    @(37, 6, 7),
    @(38, 6, 7),
    @(39, 6, 7),
    @(40, 5, 5),
    @(41, 6, 7),
    @(42, 4, 5),
    @(43, 12, 12),
    @(44, 6, 7),
    @(45, 8, 8),
    @(46, 7, 7),
    @(47, 5, 5),
    @(48, 1, 1),
    @(49, 5, 6),
    @(50, 6, 6),
    @(51, 4, 4),
    @(52, 2, 3),
    @(53, 7, 7),
    @(54, 1, 2),
    @(55, 7, 7),
    @(56, 6, 7),
    @(57, 6, 7),
    @(58, 8, 9),
    @(59, 8, 8),
    @(60, 5, 5),
    @(61, 5, 5),
    @(62, 9, 9),
    @(63, 8, 8),
    @(64, 7, 7),
    @(65, 7, 7)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value  = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
